# PlayerPerformance_4362.xlsx - additional scraping update
#
# Summary of changes:
#  1. New "Player Info" sheet inserted at the front (ID / NAME / BATTING_HAND / BOWL_STYLE).
#  2. Existing "ODI Batting" sheet: MATCH_CARD_LINK column -> MATCH_CODE (URL replaced by the
#     bare numeric match code).
#  3. Existing "ODI Bowling" sheet: same MATCH_CARD_LINK -> MATCH_CODE transformation.
#  4. New "ODI Batting Extra" sheet appended at the end with extra per-match batting stats.

$wb = $excel.ActiveWorkbook

function Set-AsText($cell, [string]$text) {
    # Writing a numeric-looking string through .Value normally auto-converts it to a
    # number (regular Excel COM behaviour). Forcing the NumberFormat to Text first keeps
    # it a string; flipping the format/style back afterwards avoids leaving a stray
    # "Text" number-format behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

function Copy-HeaderStyle($srcCell, $destCell, [string]$text) {
    $srcCell.Copy($destCell)
    $destCell.Value = $text
}

# ---------------------------------------------------------------------------
# 1) Rename MATCH_CARD_LINK -> MATCH_CODE and collapse the howstat URL down to
#    just the numeric MatchCode query parameter, on both existing sheets.
# ---------------------------------------------------------------------------

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingLinkCol = 4   # column D

$lastRow = $battingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, $battingLinkCol)
    $val = $cell.Value2
    if ($val -match "MatchCode=(\d+)") {
        Set-AsText $cell $matches[1]
    }
}
$battingSheet.Cells.Item(1, $battingLinkCol).Value = "MATCH_CODE"

# Rows that never batted in that match leave INNING_NUMBER (column B) blank;
# collapse those to genuinely-empty cells instead of an empty string.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 2)
    if ($cell.Value2 -eq "") {
        $cell.ClearContents()
    }
}

$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingLinkCol = 2   # column B

$lastRow = $bowlingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, $bowlingLinkCol)
    $val = $cell.Value2
    if ($val -match "MatchCode=(\d+)") {
        Set-AsText $cell $matches[1]
    }
}
$bowlingSheet.Cells.Item(1, $bowlingLinkCol).Value = "MATCH_CODE"

# ---------------------------------------------------------------------------
# 2) New "Player Info" sheet, inserted as the first tab.
# ---------------------------------------------------------------------------

$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$headerStyleSrc = $battingSheet.Range("A1")

$playerHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $playerHeaders.Count; $i++) {
    Copy-HeaderStyle $headerStyleSrc $playerInfo.Cells.Item(1, $i + 1) $playerHeaders[$i]
}

Set-AsText $playerInfo.Cells.Item(2, 1) "4362"
$playerInfo.Cells.Item(2, 2).Value = "Soumya Sarkar"
$playerInfo.Cells.Item(2, 3).Value = "Left Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Medium Fast"

# ---------------------------------------------------------------------------
# 3) New "ODI Batting Extra" sheet, appended as the last tab.
# ---------------------------------------------------------------------------

$lastTab = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastTab)
$extra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $extraHeaders.Count; $i++) {
    Copy-HeaderStyle $headerStyleSrc $extra.Cells.Item(1, $i + 1) $extraHeaders[$i]
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# BATTING_POSITION is a genuine number; the rest (even the purely-numeric-looking
# ones) are stored as text, matching the source scrape.
$extraRows = @(
    @("4286", "2", "9", "1", "27.65%", "NO"),
    @("4293", "2", "4", "2", "21.77%", "NO"),
    @("4296", "2", "9", "3", "30.99%", "NO"),
    @("4307", "", "", "", "", "NO"),
    @("4311", "", "", "", "", "NO"),
    @("4314", "2", "0", "0", "0.71%", "NO"),
    @("4325", "", "", "", "", "NO"),
    @("4329", "", "", "", "", "NO"),
    @("4335", "5", "0", "0", "1.15%", "NO"),
    @("4345", "2", "4", "0", "11.54%", "NO"),
    @("4349", "2", "4", "0", "9.95%", "NO"),
    @("4356", "", "", "", "", "NO"),
    @("4357", "", "", "", "", "NO"),
    @("4358", "", "", "", "", "NO"),
    @("4443", "", "", "", "", "NO"),
    @("4445", "6", "", "", "", "NO"),
    @("4447", "7", "0", "0", "2.36%", "NO"),
    @("4452", "3", "0", "0", "", "NO"),
    @("4453", "3", "3", "1", "11.81%", "NO"),
    @("4455", "3", "0", "0", "0.65%", "NO")
)

for ($i = 0; $i -lt $extraRows.Count; $i++) {
    $row = $i + 2
    $values = $extraRows[$i]

    Set-AsText $extra.Cells.Item($row, 1) $values[0]

    if ($values[1] -ne "") {
        $extra.Cells.Item($row, 2).Value = [int]$values[1]
    }

    if ($values[2] -ne "") { Set-AsText $extra.Cells.Item($row, 3) $values[2] }
    if ($values[3] -ne "") { Set-AsText $extra.Cells.Item($row, 4) $values[3] }
    if ($values[4] -ne "") { Set-AsText $extra.Cells.Item($row, 5) $values[4] }
    if ($values[5] -ne "") { Set-AsText $extra.Cells.Item($row, 6) $values[5] }
}

# ---------------------------------------------------------------------------
# Restore the original active tab (first sheet) / selection.
# ---------------------------------------------------------------------------
$playerInfo.Activate()
$playerInfo.Range("A1").Select()
